$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 ---------------------------------------------------------------
$ws.Range("A1").Value = "v"
$ws.Range("B1").Value = "ar2"
$ws.Range("C1").Value = "A/C"
$ws.Range("D1").Value = 30
$ws.Range("E1").Value = $true
# Column F is no longer used - drop its old boolean value.
$ws.Range("F1").ClearContents()

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = "v"
$ws.Range("B2").Value = "a"
$ws.Range("C2").Value = "Lâmpada"
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = $false
$ws.Range("F2").ClearContents()

# --- Row 3 (new) -----------------------------------------------------------
$ws.Range("A3").Value = "a"
$ws.Range("B3").Value = "ae"
$ws.Range("C3").Value = "A/C"
$ws.Range("D3").Value = 26
$ws.Range("E3").Value = $false
